$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metadata")

# Fill in the metadata "Value" column (column B) with the data for this
# particular indicator (URGR_013_to_031 - Urban Green).
$ws.Range("B2").Value = "URGR_013_to_031"
$ws.Range("B3").Value = "Urban Green"
$ws.Range("B4").Value = "Norway"
$ws.Range("B5").Value = "Europe"
$ws.Range("B6").Value = "B2 - Structural State Characteristics"
$ws.Range("B7").Value = "Terrestrial (T)"
$ws.Range("B8").Value = "T7 Intensive land-use biome"
$ws.Range("B9").Value = "T7.4 Urban and industrial"
$ws.Range("B10").Value = 2024
$ws.Range("B11").Value = 2025
$ws.Range("B14").Value = "First draft version, work in progress"
$ws.Range("B15").Value = "Clappe, S., Czúcz, B."
$ws.Range("B17").Value = "No"
$ws.Range("B18").Value = "National scale - Aggregated at Regional level"

# Switch the active sheet/selection back to the metadata sheet, matching
# where the author ended up after entering the data.
$ws.Activate()
$ws.Range("B21").Select()
